$wb = $excel.ActiveWorkbook

# Sheets "展览" and "全部类型" both contain the same event table in rows 2-9,
# column F ("想去人数" / number of people interested). Reset those counts to 0.
$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    for ($row = 2; $row -le 9; $row++) {
        $ws.Cells.Item($row, 6).Value = 0
    }
}
